# Apply the FFV Schedule StructureDefinition metadata refresh:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publication timestamp
#  - set Publisher value to "Alvearie Team"
#  - replace the (duplicated) Contact rows with a single Jurisdiction row
#  - update the root Extension row's Short/Definition text on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 currently holds the first "Contact" / "No display for ContactDetail" pair;
# turn it into the new Jurisdiction row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 is the duplicate "Contact" / "No display for ContactDetail" row - remove it,
# shifting the remaining rows (Description ... Context) up by one.
$meta.Rows(11).Delete()

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "FFV Schedule"
$elements.Range("L2").Value = "Schedule code for the FFV initiative"
